# Task Scheduler update: record the latest UserName / WorkGroup run into
# row 2 of Sheet1 (A2 = UserName<timestamp>, C2 = WorkGroup<timestamp>).
#
# Each scheduled run stamps a fresh unique UserName/WorkGroup pair. This
# replay performs both the prior run and the newest run in sequence so the
# worksheet ends up holding the latest pair, matching the Task Scheduler's
# "last write wins" behaviour for these two cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previous scheduled run's stamp (superseded below).
$ws.Range("A2").Value = "UserName1551090484881"
$ws.Range("C2").Value = "WorkGroup1551090662155"

# Latest scheduled run's stamp - final values for A2 / C2.
$ws.Range("A2").Value = "UserName1551346091224"
$ws.Range("C2").Value = "WorkGroup1551346328646"
